# ---------------------------------------------------------------------------
# Applies the "Update links.md, remove extra links from main doc" commit:
#   1. Bump the Date paragraph to 2017-04-04.
#   2. Add a new "Del Ponte ..." bullet right after the "Sparks et al. ..."
#      bullet in the WHAT IS THE STATE OF REPRODUCIBLE RESEARCH... list.
#   3. Remove the whole "Notes" sub-section (heading + intro paragraph + the
#      two resource hyperlinks) that used to sit between ACKOWLEDGEMENTS and
#      LITERATURE CITED.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Helper: find the 1-based Paragraphs index that contains document position $pos.
function Get-ParaIndexAt($pos) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $pr = $d.Paragraphs.Item($i).Range
        if ($pos -ge $pr.Start -and $pos -lt $pr.End) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. Date: 2016-12-07 -> 2017-04-04
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2016-12-07", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2017-04-04", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert the new "Del Ponte provides data ..." list item right after the
#    "Sparks et al. (2011, 2014) ..." bullet (same list / same Compact style).
# ---------------------------------------------------------------------------
$hyperlinks = $d.Hyperlinks
$sparksHyperlink = $null
for ($i = 1; $i -le $hyperlinks.Count; $i++) {
    $h = $hyperlinks.Item($i)
    if ($h.Address -eq "http://adamhsparks.github.io/Global-Late-Blight-MetaModelling/") {
        $sparksHyperlink = $h
        break
    }
}

$sparksParaIdx = Get-ParaIndexAt($sparksHyperlink.Range.Start)
$sparksPara = $d.Paragraphs.Item($sparksParaIdx)

$insertionPoint = $sparksPara.Range.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$delPonteIdx = $sparksParaIdx + 1
$delPontePara = $d.Paragraphs.Item($delPonteIdx)
$delPontePara.Range.Text = "Del Ponte provides data and a reproducible report that explan in details all steps of the analysis and the R codes for conducting a meta-analysis for assessing heterogeneity in relationship between white mold incidence and soybean yield and between incidence and soybean yied."

# ---------------------------------------------------------------------------
# 3. Delete the "Notes" sub-section: from the "Notes" Heading3 paragraph up
#    to (but not including) the "LITERATURE CITED" Heading3 paragraph.
# ---------------------------------------------------------------------------
$notesBookmark = $d.Bookmarks.Item("notes")
$literatureBookmark = $d.Bookmarks.Item("literature-cited")

$notesParaIdx = Get-ParaIndexAt($notesBookmark.Start)
$literatureParaIdx = Get-ParaIndexAt($literatureBookmark.Start)

$deleteStart = $d.Paragraphs.Item($notesParaIdx).Range.Start
$deleteEnd = $d.Paragraphs.Item($literatureParaIdx).Range.Start
$d.Range($deleteStart, $deleteEnd).Delete()

Write-Output "Edit complete."
